$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.082804203033447
$ws.Range("B1").Value = 3.550858497619629
$ws.Range("C1").Value = 3.149245738983154
$ws.Range("D1").Value = 3.388283967971802
$ws.Range("E1").Value = 1.020687460899353
